$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("331:332").Insert()

$ws.Range("A331").Value = 5
$ws.Range("B331").Value = "Macroferia Regional de Talca"
$ws.Range("C331").Value = "Maule"
$ws.Range("D331").Value = 44461
$ws.Range("E331").Value = 7
$ws.Range("F331").Value = "Fruta"
$ws.Range("G331").Value = 100108
$ws.Range("H331").Value = "Tropicales y subtropicales"
$ws.Range("I331").Value = 100108006
$ws.Range("J331").Value = "Plátano"
$ws.Range("K331").Value = "Sin especificar"
$ws.Range("L331").Value = "Pintón"
$ws.Range("M331").Value = 300
$ws.Range("N331").Value = 14000
$ws.Range("O331").Value = 14000
$ws.Range("P331").Value = 14000
$ws.Range("Q331").Value = "$/caja 20 kilos"
$ws.Range("R331").Value = "Ecuador"
$ws.Range("S331").Value = 700
$ws.Range("T331").Value = 20

$ws.Range("A332").Value = 5
$ws.Range("B332").Value = "Macroferia Regional de Talca"
$ws.Range("C332").Value = "Maule"
$ws.Range("D332").Value = 44461
$ws.Range("E332").Value = 7
$ws.Range("F332").Value = "Fruta"
$ws.Range("G332").Value = 100108
$ws.Range("H332").Value = "Tropicales y subtropicales"
$ws.Range("I332").Value = 100108006
$ws.Range("J332").Value = "Plátano"
$ws.Range("K332").Value = "Sin especificar"
$ws.Range("L332").Value = "Primera Pintón"
$ws.Range("M332").Value = 550
$ws.Range("N332").Value = 15000
$ws.Range("O332").Value = 15000
$ws.Range("P332").Value = 15000
$ws.Range("Q332").Value = "$/caja 20 kilos"
$ws.Range("R332").Value = "Ecuador"
$ws.Range("S332").Value = 750
$ws.Range("T332").Value = 20
